# Timesheet changes by Ruchika(MT2012119)
# Mark columns AC:AL and AO as "OFF" (reusing the existing OFF formatting
# already present on cells such as Q28) for rows 28-31, and fill in the
# AM/AN totals for those same rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$offRows = 28, 29, 30, 31

# Copy the formatting already used for "OFF" cells (e.g. Q28, which is
# styled with the grey OFF fill) onto the new OFF cells, then stamp in
# the "OFF" text.
$ws.Range("Q28").Copy() | Out-Null
foreach ($r in $offRows) {
    $ws.Range("AC$r`:AL$r").PasteSpecial(-4122)
    $ws.Range("AO$r").PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

foreach ($r in $offRows) {
    $ws.Range("AC$r`:AL$r").Value = "OFF"
    $ws.Range("AO$r").Value = "OFF"
}

# Fill in the AM/AN values for each row (style stays as-is, s=26).
$ws.Range("AM28").Value = 0
$ws.Range("AN28").Value = 0

$ws.Range("AM29").Value = 1
$ws.Range("AN29").Value = 1

$ws.Range("AM30").Value = 0
$ws.Range("AN30").Value = 0

$ws.Range("AM31").Value = 0
$ws.Range("AN31").Value = 0
